# CasosColombia.xlsx update:
#  1) Replace 7 scattered numeric "2"/"3"/"30"/"32"/"35"/"206" outlier cells
#     with the text value "NaN" (these columns already use "NaN" elsewhere
#     in the sheet to flag missing/unreliable daily counts).
#  2) Append a new data row (row 192, date 2020-09-12) with the day's case
#     counts for every one of the 128 tracked columns (A:DX).
#  3) Leave the freeze-pane split as-is (first row/column) and move the
#     active selection down to the newly typed-in last cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Flip the six/seven outlier numeric cells over to the "NaN" marker text ---
$naCells = @("BZ18", "CB20", "CF93", "AI105", "AI106", "AE130", "AF157")
foreach ($cellRef in $naCells) {
    $ws.Range($cellRef).Value = "NaN"
}

# --- 2) Append row 192 across all 128 columns (A:DX) ---
$row192 = @(
    44086,708964,2725,94027,65726,238334,26842,4912,3884,7249,
    7186,15319,3859,22024,27987,6510,7388,13779,11360,15997,
    13525,3348,1969,8004,24166,13078,9415,53033,1509,482,
    578,459,342,223,457,2003,4251,37026,8166,2505,
    41396,1035,21584,1487,9311,1596,1586,6227,1765,955,
    2484,2652,54399,13310,4313,8641,5319,281,1435,2649,
    736,2116,9160,9125,9707,14084,1925,860,11399,9072,
    10572,2013,1815,4396,4131,1403,5206,2945,1711,838,
    2549,2128,1625,1268,5874,1807,1313,1557,1903,1862,
    2200,1370,1164,1156,755,3183,1258,877,900,1621,
    1424,709,814,1120,1416,1227,1337,1048,326,349,
    762,684,451,535,363,647,729,519,484,372,
    518,128176,301117,14130,129119,80075,37740,10939
)

$newRow = 192
for ($i = 0; $i -lt $row192.Length; $i++) {
    $ws.Cells.Item($newRow, $i + 1).Value = $row192[$i]
}

# --- 3) Move the active selection to the last cell that was filled in ---
$ws.Range("DX192").Select()
